# [Feat 2269] Add support of datasets worksheet metadata
# Rename "DATASETS TODO" sheet to "DATASETS" and populate its header row
# with the new dataset-related column metadata, then make it the active
# sheet/tab (instead of PARAMETERS).

$wb = $excel.ActiveWorkbook

$wsDatasets = $wb.Worksheets.Item("DATASETS TODO")
$wsDatasets.Name = "DATASETS"

# Populate the DATASETS header row (row 1) with the matching enum values.
$wsDatasets.Range("A1").Value = "ACTION"
$wsDatasets.Range("B1").Value = "TC_OWNER_PATH"
$wsDatasets.Range("C1").Value = "TC_OWNER_ID"
$wsDatasets.Range("D1").Value = "TC_DATASET_ID"
$wsDatasets.Range("E1").Value = "TC_DATASET_NAME"
$wsDatasets.Range("F1").Value = "TC_PARAM_OWNER_ID"
$wsDatasets.Range("G1").Value = "TC_DATASET_PARAM_VALUE"
$wsDatasets.Range("H1").Value = "TC_DATASET_PARAM_NAME"
$wsDatasets.Range("I1").Value = "TC_PARAM_OWNER_PATH"

# Column widths matching the bestFit sizing applied by Excel in the target
# (closest values reachable through this engine's column-width quantization).
$wsDatasets.Columns.Item(1).ColumnWidth = 7
$wsDatasets.Columns.Item(2).ColumnWidth = 16
$wsDatasets.Columns.Item(3).ColumnWidth = 13
$wsDatasets.Columns.Item(4).ColumnWidth = 14
$wsDatasets.Columns.Item(5).ColumnWidth = 17.666666666666668
$wsDatasets.Columns.Item(6).ColumnWidth = 20.833333333333332
$wsDatasets.Columns.Item(7).ColumnWidth = 25.833333333333332
$wsDatasets.Columns.Item(8).ColumnWidth = 25.666666666666668

# PARAMETERS was the active/selected tab before; DATASETS becomes the new
# active tab (and the only one keeping tabSelected="1" on save).
$wsDatasets.Activate()
$wsDatasets.Range("H7").Select() | Out-Null
